# Finishing touches on HIV progression
# The "constants" sheet contained four comorbidity/HIV related rows:
#   row 49: comorb_startage_hiv                     = 10
#   row 50: comorb_multiplier_diabetes_progression   = 3.11
#   row 51: comorb_multiplier_hiv_progression        = 3.67
#   row 52: comorb_multiplier_hiv_late_progression   = 36.7
#
# The edit drops the now-unused "comorb_startage_hiv" parameter and the old
# "comorb_multiplier_hiv_progression" value (3.67), keeping the former "late
# progression" row (value 36.7, along with its number formatting) but
# renaming it to just "comorb_multiplier_hiv_progression", leaving:
#   row 49: comorb_multiplier_diabetes_progression   = 3.11
#   row 50: comorb_multiplier_hiv_progression        = 36.7
# and every row below shifts up by two.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("constants")

# Delete row 51 first (old comorb_multiplier_hiv_progression, value 3.67) so
# that row indices above it (including row 49) are unaffected by this
# deletion.
$ws.Rows.Item(51).Delete()

# Delete row 49 (comorb_startage_hiv, value 10). After this, the old row 52
# (comorb_multiplier_hiv_late_progression, value 36.7) becomes row 50.
$ws.Rows.Item(49).Delete()

# The row that used to be "comorb_multiplier_hiv_late_progression" now sits
# at row 50 with its original value/formatting intact; rename its label so
# it reads as the (now sole) "comorb_multiplier_hiv_progression" parameter.
$ws.Range("A50").Value = "comorb_multiplier_hiv_progression"

# Match the selection left behind by the editor.
$ws.Activate()
$ws.Range("B48").Select()
